# Edit the "Worksheet" sheet: fill in a handful of previously-blank / placeholder
# cells in column C. Cells that look numeric ("32", "-8", ...) must be forced to
# Text so they keep matching the sheet's existing "numbers stored as text" data,
# without disturbing the cell's current style index (NumberFormat is flipped to
# Text just for the assignment, then restored to General).
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Worksheet")

function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.NumberFormat = "General"
}

Set-TextValue $ws.Range("C14") "32"
Set-TextValue $ws.Range("C15") "/"

Set-TextValue $ws.Range("C19") "-8"
Set-TextValue $ws.Range("C20") "/"

Set-TextValue $ws.Range("C24") "-8"
Set-TextValue $ws.Range("C25") "/"

Set-TextValue $ws.Range("C29") "-8"
Set-TextValue $ws.Range("C30") "/"

Set-TextValue $ws.Range("C34") "-8"
Set-TextValue $ws.Range("C35") "/"

Set-TextValue $ws.Range("C39") "-8"
Set-TextValue $ws.Range("C40") "/"

Set-TextValue $ws.Range("C44") "-8"
Set-TextValue $ws.Range("C45") "/"

Set-TextValue $ws.Range("C49") "-8"
Set-TextValue $ws.Range("C50") "/"

Set-TextValue $ws.Range("C54") "-8"
Set-TextValue $ws.Range("C55") "/"

Set-TextValue $ws.Range("C59") "-8"
Set-TextValue $ws.Range("C60") "/"

# Append a duplicate of "Worksheet" (with the edits above baked in) as a new
# last sheet named "test1", matching the new sheetId/rId appended to the
# workbook.
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws.Copy($null, $lastSheet)
$newSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet.Name = "test1"
